$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value would otherwise be auto-parsed as a number by Excel;
# temporarily force text format so the literal string is preserved, then clear the
# number format again afterwards so the cell style matches the untouched cells.
$textForced = @("D5", "D6", "D13", "D14", "D21", "D22", "D23", "D25", "D27", "D28", "D32", "D33", "D34", "D38", "D40", "D43", "D44", "D45", "D46", "D48", "D50")
foreach ($addr in $textForced) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "69.222.84"
$ws.Range("E2").Value = "  +0.48%  "
$ws.Range("D3").Value = "3.933.96"
$ws.Range("E3").Value = "  +5.20%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "604.62"
$ws.Range("E5").Value = "  +0.69%  "
$ws.Range("D6").Value = "164.67"
$ws.Range("E6").Value = "  -0.43%  "
$ws.Range("D7").Value = "3.933.08"
$ws.Range("E7").Value = "  +5.21%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  -1.11%  "
$ws.Range("E10").Value = "  -1.32%  "
$ws.Range("E11").Value = "  -0.34%  "
$ws.Range("E12").Value = "  +1.02%  "
$ws.Range("D13").Value = "37.19"
$ws.Range("E13").Value = "  -1.37%  "
$ws.Range("D14").Value = "0.0000247"
$ws.Range("E14").Value = "  -0.33%  "
$ws.Range("D15").Value = "4.582.01"
$ws.Range("E15").Value = "  +4.97%  "
$ws.Range("D16").Value = "3.927.26"
$ws.Range("E16").Value = "  +5.03%  "
$ws.Range("D17").Value = "69.256.19"
$ws.Range("E17").Value = "  +0.38%  "
$ws.Range("E18").Value = "  +1.22%  "
$ws.Range("E19").Value = "  -0.22%  "
$ws.Range("E20").Value = "  -1.80%  "
$ws.Range("D21").Value = "11.28"
$ws.Range("E21").Value = "  +1.46%  "
$ws.Range("D22").Value = "489.62"
$ws.Range("E22").Value = "  -0.52%  "
$ws.Range("D23").Value = "0.726"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("E24").Value = "  +13.88%  "
$ws.Range("D25").Value = "84.57"
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("E26").Value = "  +0.59%  "
$ws.Range("D27").Value = "12.17"
$ws.Range("E27").Value = "  -1.13%  "
$ws.Range("D28").Value = "10.14"
$ws.Range("E28").Value = "  +1.17%  "
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("E30").Value = "  -0.20%  "
$ws.Range("D31").Value = "4.081.60"
$ws.Range("E31").Value = "  +4.91%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "32.65"
$ws.Range("E32").Value = "  +3.55%  "
$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").Value = "7.91"
$ws.Range("E33").Value = "  -3.61%  "
$ws.Range("D34").Value = "2.40"
$ws.Range("E34").Value = "  -1.45%  "
$ws.Range("D35").Value = "3.871.48"
$ws.Range("E35").Value = "  +5.33%  "
$ws.Range("E36").Value = "  -0.22%  "
$ws.Range("E37").Value = "  +3.71%  "
$ws.Range("D38").Value = "0.140"
$ws.Range("E38").Value = "  +1.88%  "
$ws.Range("E39").Value = "  +0.39%  "
$ws.Range("D40").Value = "0.998"
$ws.Range("E40").Value = "  -0.25%  "
$ws.Range("E41").Value = "  -0.60%  "
$ws.Range("E42").Value = "  -3.17%  "
$ws.Range("D43").Value = "441.09"
$ws.Range("E43").Value = "  +4.03%  "
$ws.Range("D44").Value = "2.01"
$ws.Range("E44").Value = "  +1.04%  "
$ws.Range("D45").Value = "48.48"
$ws.Range("E45").Value = "  -0.16%  "
$ws.Range("D46").Value = "8.48"
$ws.Range("E46").Value = "  +0.58%  "
$ws.Range("D48").Value = "27.79"
$ws.Range("E48").Value = "  +19.57%  "
$ws.Range("D49").Value = "2.850.38"
$ws.Range("E49").Value = "  +2.64%  "
$ws.Range("D50").Value = "141.73"
$ws.Range("E50").Value = "  +0.51%  "
$ws.Range("E51").Value = "  +2.40%  "

foreach ($addr in $textForced) {
    $ws.Range($addr).Style = "Normal"
}
